# SARAALERT-1521: Import/export follow-up fields
# Adds two new trailing columns to the Sara Alert "Invalid Monitorees"
# import/export template: "Follow-Up Reason" and "Follow-Up Note".
# These land immediately after the existing "Vaccine 2 Notes" column
# (the last used column, DH), becoming DI and DJ respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1) - appended right after the last existing
# header ("Vaccine 2 Notes" in DH1).
$ws.Range("DI1").Value = "Follow-Up Reason"
$ws.Range("DJ1").Value = "Follow-Up Note"

# Match the bestFit/custom widths the real workbook ends up with for
# these two new columns (computed from the header text, same as the
# neighbouring bestFit columns).
$ws.Columns.Item(113).ColumnWidth = 13.998697916666666
$ws.Columns.Item(114).ColumnWidth = 12.166666666666666
